$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for the new rows from the most similar existing rows ---
# Rows 39-40 should look like the existing "green" rows (e.g. row 33)
$ws.Range("A33:D33").Copy()
$ws.Range("A39:D40").PasteSpecial(-4122)

# Row 41 should look like the existing "yellow" rows (e.g. row 38)
$ws.Range("A38:D38").Copy()
$ws.Range("A41:D41").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 39 ---
$ws.Range("A39").Value = "Filtragem incorreta ao buscar os despachos de patentes no cadastro do processo"
$ws.Range("B39").Value = "Defeito"
$ws.Range("C39").Value = "Desenvolvido"
$ws.Range("D39").Value = "N/A"

# --- Row 40 ---
$ws.Range("A40").Value = "Filtragem incorreta das pastas cadastradas"
$ws.Range("B40").Value = "Defeito"
$ws.Range("C40").Value = "Desenvolvido"
$ws.Range("D40").Value = "N/A"

# --- Row 41 ---
$ws.Range("A41").Value = "Implementado novo tratamento do número do processo na leitura da revista"
$ws.Range("B41").Value = "Melhoria"
$ws.Range("C41").Value = "Desenvolvido"
$ws.Range("D41").Value = "N/A"

# --- Update the visible selection to match the new bottom of the list ---
$ws.Range("A42").Select()
